# Add data for 2022-10-18 (extend "through October 09" -> "through October 10")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (tab name) to reflect the new date.
$ws.Name = "Through 2022-10-10"

# 2. Update the header text in B1 (shared string) to match the new date.
$ws.Range("B1").Value = "October 2022 (through October 10)"

# 3. New incidents added on 2022-10-10 in column B (current/in-progress month)
#    for the affected neighborhoods (rows 2, 5, 7, 32, 35).
$ws.Range("B2").Value = 1
$ws.Range("B5").Value = 1
$ws.Range("B7").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("B35").Value = 1

# 4. Increment existing monthly totals that also changed by +1 (these are
#    historical-month columns whose counts were revised upward).
$ws.Range("L2").Value = 6
$ws.Range("L3").Value = 3
$ws.Range("AP3").Value = 2
$ws.Range("L4").Value = 4
$ws.Range("V5").Value = 6
$ws.Range("L6").Value = 5
$ws.Range("BJ7").Value = 2
$ws.Range("BT7").Value = 2
$ws.Range("V12").Value = 1
$ws.Range("V14").Value = 1
$ws.Range("V16").Value = 1
$ws.Range("AF17").Value = 1
$ws.Range("AF23").Value = 1
$ws.Range("AF25").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("L28").Value = 3
$ws.Range("AP29").Value = 2
